$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-23 11:07:19"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA" + $row).Value = $newTimestamp
    }
}
